$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 19828.5
$ws.Range("J3").Value = 19828.5
$ws.Range("L3").Value = 19828.5
$ws.Range("N3").Value = -20056.5
$ws.Range("H55").Value = 576.7143
$ws.Range("I55").Value = 372.5
$ws.Range("J55").Value = 849
$ws.Range("K55").Value = 372.5
$ws.Range("L55").Value = 849
$ws.Range("M55").Value = -158.5
$ws.Range("N55").Value = -1277
$ws.Range("H100").Value = 2360.5
$ws.Range("I100").Value = 2230.6667
$ws.Range("K100").Value = 2230.6667
$ws.Range("M100").Value = -1689.6667
$ws.Range("H102").Value = 19828.5
$ws.Range("J102").Value = 19828.5
$ws.Range("L102").Value = 19828.5
$ws.Range("N102").Value = -26318.5
$ws.Range("H129").Value = 1543.2222
$ws.Range("I129").Value = 815.1667
$ws.Range("K129").Value = 2445.5001
$ws.Range("M129").Value = 2554.4999
$ws.Range("H137").Value = 2128.0952
$ws.Range("I137").Value = 1510
$ws.Range("K137").Value = 4530
$ws.Range("M137").Value = -1980

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 4225.75
$ws.Range("J3").Value = 4225.75
$ws.Range("L3").Value = 4225.75
$ws.Range("N3").Value = -4455.75
$ws.Range("H32").Value = 1652.9048
$ws.Range("I32").Value = 1514.4324
$ws.Range("K32").Value = 1514.4324
$ws.Range("M32").Value = -1227.4324
$ws.Range("H33").Value = 16513
$ws.Range("I33").Value = 26
$ws.Range("K33").Value = 26
$ws.Range("M33").Value = 303
$ws.Range("H36").Value = 50000
$ws.Range("I36").Value = 50000
$ws.Range("K36").Value = 50000
$ws.Range("M36").Value = -49654
$ws.Range("H45").Value = 1603
$ws.Range("I45").Value = 1637.6666
$ws.Range("J45").Value = 1499
$ws.Range("K45").Value = 1637.6666
$ws.Range("L45").Value = 1499
$ws.Range("M45").Value = -1260.6666
$ws.Range("N45").Value = -2253
$ws.Range("H74").Value = 940.1
$ws.Range("I74").Value = 959.5
$ws.Range("K74").Value = 959.5
$ws.Range("M74").Value = -85.5
$ws.Range("H77").Value = 940.1
$ws.Range("I77").Value = 959.5
$ws.Range("K77").Value = 4797.5
$ws.Range("M77").Value = -429.5
$ws.Range("H88").Value = 2667.8333
$ws.Range("I88").Value = 220
$ws.Range("J88").Value = 3157.4
$ws.Range("K88").Value = 220
$ws.Range("L88").Value = 3157.4
$ws.Range("M88").Value = 186
$ws.Range("N88").Value = -3969.4
$ws.Range("H91").Value = 2667.8333
$ws.Range("I91").Value = 220
$ws.Range("J91").Value = 3157.4
$ws.Range("K91").Value = 220
$ws.Range("L91").Value = 3157.4
$ws.Range("M91").Value = 1184
$ws.Range("N91").Value = -5965.4
$ws.Range("H95").Value = 15000
$ws.Range("J95").Value = 15000
$ws.Range("L95").Value = 15000
$ws.Range("N95").Value = -20492
$ws.Range("H96").Value = 21333.334
$ws.Range("J96").Value = 21333.334
$ws.Range("L96").Value = 21333.334
$ws.Range("N96").Value = -26825.334
$ws.Range("H122").Value = 3383.6667
$ws.Range("I122").Value = 1935.3334
$ws.Range("K122").Value = 5806.0002
$ws.Range("M122").Value = -3356.0002
$ws.Range("H132").Value = 2457.5
$ws.Range("I132").Value = 2457.5
$ws.Range("K132").Value = 7372.5
$ws.Range("M132").Value = -4842.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 2167.5
$ws.Range("I5").Value = 1700
$ws.Range("J5").Value = 2635
$ws.Range("K5").Value = 1700
$ws.Range("L5").Value = 2635
$ws.Range("M5").Value = -1587
$ws.Range("N5").Value = -2861
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").ClearContents()
$ws.Range("N40").Value = 0
$ws.Range("H96").Value = 19950
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("H101").Value = 0
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 0
$ws.Range("L101").ClearContents()
$ws.Range("M101").ClearContents()
$ws.Range("N101").Value = 0
$ws.Range("H132").Value = 200640
$ws.Range("J132").Value = 200640
$ws.Range("L132").Value = 200640
$ws.Range("N132").Value = -210760

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 1894.2
$ws.Range("J3").Value = 1894.2
$ws.Range("L3").Value = 1894.2
$ws.Range("N3").Value = -2120.2
$ws.Range("H17").Value = 3243.6
$ws.Range("I17").Value = 2054.4443
$ws.Range("K17").Value = 2054.4443
$ws.Range("M17").Value = -1880.4443
$ws.Range("H88").Value = 57996.668
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 57996.668
$ws.Range("K88").Value = 0
$ws.Range("L88").ClearContents()
$ws.Range("M88").Value = 57996.668
$ws.Range("N88").Value = -58808.668
$ws.Range("H91").Value = 57996.668
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 57996.668
$ws.Range("K91").Value = 0
$ws.Range("L91").ClearContents()
$ws.Range("M91").Value = 57996.668
$ws.Range("N91").Value = -60804.668
$ws.Range("H99").Value = 5331.75
$ws.Range("I99").Value = 4824.75
$ws.Range("J99").Value = 5838.75
$ws.Range("K99").Value = 4824.75
$ws.Range("L99").Value = 5838.75
$ws.Range("M99").Value = -3326.75
$ws.Range("N99").Value = -8834.75
$ws.Range("H126").Value = 5331.75
$ws.Range("I126").Value = 4824.75
$ws.Range("J126").Value = 5838.75
$ws.Range("K126").Value = 14474.25
$ws.Range("L126").Value = 17516.25
$ws.Range("M126").Value = -12004.25
$ws.Range("N126").Value = -22456.25

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 112636.63
$ws.Range("J4").Value = 2232.889
$ws.Range("L4").Value = 6698.667
$ws.Range("N4").Value = -6922.667
$ws.Range("H109").Value = 750.5
$ws.Range("I109").Value = 750.5
$ws.Range("K109").Value = 2251.5
$ws.Range("M109").Value = -1211.5
$ws.Range("H131").Value = 0
$ws.Range("I131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("M131").ClearContents()
$ws.Range("H134").Value = 168622
$ws.Range("I134").Value = 200339.8
$ws.Range("J134").Value = 10033
$ws.Range("K134").Value = 601019.3999999999
$ws.Range("L134").Value = 30099
$ws.Range("M134").Value = -595949.3999999999
$ws.Range("N134").Value = -40239
$ws.Range("H139").Value = 1536.75
$ws.Range("I139").Value = 1132.3334
$ws.Range("K139").Value = 3397.0002
$ws.Range("M139").Value = 1742.9998
$ws.Range("H140").Value = 1691.6666
$ws.Range("I140").Value = 1691.6666
$ws.Range("K140").Value = 5074.9998
$ws.Range("M140").Value = 105.0002000000004
$ws.Range("H141").Value = 2650
$ws.Range("I141").Value = 2650
$ws.Range("K141").Value = 7950
$ws.Range("M141").Value = -2770

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 951.5
$ws.Range("I4").Value = 903
$ws.Range("K4").Value = 903
$ws.Range("M4").Value = -791
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").ClearContents()
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = 0
$ws.Range("H132").Value = 2000
$ws.Range("I132").Value = 2000
$ws.Range("K132").Value = 6000
$ws.Range("M132").Value = -3470

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3297.3333
$ws.Range("I61").Value = 3918.8
$ws.Range("J61").Value = 190
$ws.Range("K61").Value = 3918.8
$ws.Range("L61").Value = 190
$ws.Range("M61").Value = -3716.8
$ws.Range("N61").Value = -594
$ws.Range("H82").Value = 399.5
$ws.Range("J82").Value = 399.5
$ws.Range("L82").Value = 399.5
$ws.Range("N82").Value = -1121.5
$ws.Range("H85").Value = 399.5
$ws.Range("J85").Value = 399.5
$ws.Range("L85").Value = 399.5
$ws.Range("N85").Value = -2895.5
$ws.Range("H113").Value = 3297.3333
$ws.Range("I113").Value = 3918.8
$ws.Range("J113").Value = 190
$ws.Range("K113").Value = 3918.8
$ws.Range("L113").Value = 190
$ws.Range("M113").Value = -1748.8
$ws.Range("N113").Value = -4530

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").Value = 0
